# Update "想去人数" (column F) values on the "展览" and "全部类型" sheets
# to reflect newly generated counts (gh-pages output at 456a3b4).

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 546
$ws1.Range("F10").Value = 15984
$ws1.Range("F14").Value = 6254
$ws1.Range("F31").Value = 11181

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 546
$ws4.Range("F10").Value = 15984
$ws4.Range("F14").Value = 6254
$ws4.Range("F32").Value = 11181
